$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2:D3").Value = -0.0446
$ws.Range("E2:E3").Value = 0.00488
$ws.Range("F2:F3").Value = 0.04650000000000001
$ws.Range("G2:G3").Value = 0.2066281227885116
$ws.Range("H2:H3").Value = 0.2066281227885116
$ws.Range("I2:I3").Value = 0.1567836780473476
$ws.Range("J2:J3").Value = 0.1234121871037787
$ws.Range("K2:K3").Value = 748.4
$ws.Range("L2:L3").Value = 0.09992923237151669
$ws.Range("M2:M3").Value = 411.9
$ws.Range("N2:N3").Value = 0.07403745910773989
$ws.Range("O2:O3").Value = 0.5503741314804917
$ws.Range("P2:P3").Value = 110.3
$ws.Range("Q2:Q3").Value = 0.01982600567997987
$ws.Range("R2:R3").Value = 0.147381079636558
$ws.Range("S2:S3").Value = 301.6
$ws.Range("T2:T3").Value = 0.7322165574168488
$ws.Range("U2:U3").Value = 3178
$ws.Range("V2:V3").Value = 0.5712334184132006
$ws.Range("W2:W3").Value = 0.1132001270551934
$ws.Range("X2:X3").Value = 0.08969388115417745
$ws.Range("Y2:Y3").Value = 0.02350624590101592
$ws.Range("Z2:Z3").Value = 0.5590865657380035
$ws.Range("AA2:AA3").Value = 0.06899809585806754
$ws.Range("AB2:AB3").Value = 0.0386671582247967
$ws.Range("AC2:AC3").Value = 0.03033093763327084
$ws.Range("AD2:AD3").Value = 11698.2
$ws.Range("AE2:AE3").Value = 0
$ws.Range("AF2:AF3").Value = 11698.2
$ws.Range("AG2:AG3").Value = 8520.200000000001
$ws.Range("AH2:AH3").Value = 0.6777007925105437
$ws.Range("AI2:AI3").Value = 0.6223075736377613
$ws.Range("AJ2:AJ3").Value = 0.6049731602715215
$ws.Range("AK2:AK3").Value = 0.5454638574657013
$ws.Range("AL2:AL3").Value = 323
$ws.Range("AM2:AM3").Value = 323
$ws.Range("AN2:AN3").Value = 9.738761238761239
$ws.Range("AO2:AO3").Value = 3.635294117647059
$ws.Range("AP2:AP3").Value = 7.093073593073593
$ws.Range("AQ2:AQ3").Value = 3.635294117647059